# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
# Both sheets carry the same event rows, so the same updates are applied twice.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7194
    4  = 5327
    6  = 166
    11 = 98
    12 = 196
    13 = 639
    14 = 230
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
